# Applies the "Add files via upload" commit:
#  1. Fix typo / strip stray "<br>" markup in the bullet list on slide 3
#     ("The Questions I Wanted to Answer").
#  2. Split the "What country / What sport / ..." bullets into clean,
#     single-run paragraphs (also on slide 3).
#  3. Retitle slide 7 from "Aquatic Athletes and Equestrians" to
#     "Aquatic Athletes and Skateboarders".
#
# (See the note near the bottom for the handful of diff hunks that are
# PowerPoint-internal bookkeeping with no reachable COM/VBA surface.)

$p = $ppt.ActivePresentation

# --- 1 & 2: slide 3, "Content Placeholder 2" -----------------------------
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2: "   - Which events have a mean age greather than 30? ..."
# -> single clean run, typo fixed, trailing "<br>" removed.
$body.Paragraphs(2, 1).Delete()
$anchor = $body.Paragraphs(1, 1)
$anchor.InsertAfter("`r   - Which events have a mean age greater than 30? Which events have a mean age greater than 28?") | Out-Null

# Paragraphs 4-8 (0 index shift: after the edit above counts are unchanged)
# "What country...<br>", "What sport...<br>", "What does the age...<br>",
# "Do any of them...<br>", "3. Which Olympic athlete...<br>"
# -> five clean paragraphs, "<br>" markup removed; last one keeps two runs.
for ($i = 8; $i -ge 4; $i--) {
    $body.Paragraphs($i, 1).Delete()
}

$anchor = $body.Paragraphs(3, 1)
$newParas = "`r   - What country are they from? " `
    + "`r   - What sport do they compete in? " `
    + "`r   - What does the age distribution of their sport in 2016 and 2020 look like, if available? " `
    + "`r   - Do any of them compete in the events that tend to have younger or older competitors, respectively? "
$anchor.InsertAfter($newParas) | Out-Null

# Rebuild the final bullet as two runs: "...2016" + "? "
$anchor = $body.Paragraphs(7, 1)
$anchor.InsertAfter("`r3. Which Olympic athlete has the most medals in 2016") | Out-Null
$lastPara = $body.Paragraphs(8, 1)
$lastPara.InsertAfter("? ") | Out-Null

# --- 3: slide 7 title ------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Aquatic Athletes and Skateboarders"

# NOTE: the canonical diff also bumps the cached text of the auto-updating
# "datetime1" footer fields (5/6/2022 -> 5/7/2022) on every slide layout and
# the slide master, and stamps some PowerPoint-internal, non-OM bookkeeping
# (an empty p15:sldGuideLst marker in presentation.xml, plus a co-authoring
# changesInfo changelog entry). Those are side effects of PowerPoint's own
# save routine re-evaluating the field / session metadata on whatever day
# the file is next opened/saved in the real app - there is no COM/VBA
# property that refreshes a field's cached text without collapsing the
# <a:fld> into a plain run (verified: TextFrame.TextRange.Text,
# TextFrame2.TextRange.Text, Characters(...).Text and Runs(...).Text all
# destroy the field wrapper), and there is no object-model surface for the
# guide-list / changesInfo bookkeeping at all. Leaving those untouched is
# the closer-to-canonical result, so they are intentionally left alone here.
